$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MCQ010 ---
$ws.Range("B2").Value = 0.046

# --- Row 3: MCQ160A -> MCQ195 ---
$ws.Range("C3").Value = "If MCQ195=1 use 0.089, other use given value "
$ws.Range("A3").Value = "MCQ195"
$ws.Range("B3").Value = 0.338

# --- Row 4: MCQ160N ---
$ws.Range("B4").Value = 0.438

# --- Row 5: MCQ160B ---
$ws.Range("B5").Value = 0.085

# --- Row 6: MCQ160C (was NA, now has weight + note) ---
$ws.Range("B6").Value = 0.085
$ws.Range("C6").Value = "1=yes"

# --- Row 7: MCQ160D ---
$ws.Range("B7").Value = 0.093

# --- Row 8: MCQ160E ---
$ws.Range("B8").Value = 0.432

# --- Row 9: MCQ160F (note shortened to plain "1=yes") ---
$ws.Range("B9").Value = 0.309
$ws.Range("C9").Value = "1=yes"

# --- Row 10: MCQ160M (was NA, now has weight + note) ---
$ws.Range("B10").Value = 0.082
$ws.Range("C10").Value = "1=yes"

# --- Rows 11-12: MCQ160G / MCQ160K stay NA, unchanged ---

# --- Row 13: MCQ160O ---
$ws.Range("B13").Value = 0.392

# --- Row 14: MCQ160L -> MCQ053 (was NA, now has weight + note) ---
$ws.Range("A14").Value = "MCQ053"
$ws.Range("B14").Value = 0.024
$ws.Range("C14").Value = "1=yes"

# --- Row 15: MCQ500 removed entirely ---
$ws.Range("A15:C15").ClearContents()

# --- Rows 16-17: MCQ510A / MCQ510B unchanged ---

# --- Row 18: MCQ510C ---
$ws.Range("B18").Value = 0.22

# --- Row 19: MCQ510D ---
$ws.Range("B19").Value = 0.22

# --- Row 20: MCQ510E ---
$ws.Range("B20").Value = 0.22

# --- Row 21: MCQ510f ---
$ws.Range("B21").Value = 0.22

# --- Row 22: MCQ520 (was NA, now has weight + note) ---
$ws.Range("B22").Value = 0.148
$ws.Range("C22").Value = "1=yes"

# --- Row 23: MCQ550 stays NA, unchanged ---

# --- Row 24: MCQ560 (was NA, now has weight + note) ---
$ws.Range("B24").Value = 0.15
$ws.Range("C24").Value = "1=yes"

# --- Row 25: MCQ203 (was NA, now has weight + note) ---
$ws.Range("B25").Value = 0.387
$ws.Range("C25").Value = "1=yes"

# --- Row 26: MCQ230A -> MCQ230A/B/C, merged note, B cleared ---
$newNote = "if 10: 0.293; if 11: 0.332; if 12: 0.332; if 13: 0.332; `nif 14: 0.332; if 15: 0.332; if 16: 0.26; if 17: 0.332; if 18: 0.332; if 19: 0.332; if 20: 0.332; if 21: 0.332; if 22: 0.332; if 23: 0.332; if 24: 0.332; if 25: 0.332; if 26: 0.332; if 27: 0.332; if 28: 0.332; if 29: 0.332; if 30: 0.332; if 31: 0.332; if 32: 0.332; if 33: 0.332; if 34: 0.332; if 35: 0.332; if 36: 0.332; if 37: 0.332; if 38: 0.332; if 39: 0.332"
$ws.Range("C26").Value = $newNote
$ws.Range("A26").Value = "MCQ230A/B/C"
$ws.Rows(26).RowHeight = 136

# --- Rows 27-28: MCQ230B / MCQ230C removed entirely (folded into row 26) ---
$ws.Range("A27:C28").ClearContents()

# --- View: selection moves to D26, scroll resets to top-left default ---
$ws.Range("D26").Select()
